$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.787.60'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.216.04'
$ws.Range("E3").Value = '  -2.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.04'
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.55'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.215.69'
$ws.Range("E8").Value = '  -2.31%  '
$ws.Range("E9").Value = '  +0.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.144'
$ws.Range("E10").Value = '  -3.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.30'
$ws.Range("E11").Value = '  -3.33%  '
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.07'
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.741.79'
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("E16").Value = '  -1.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.208.50'
$ws.Range("E17").Value = '  -2.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.745.95'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.68'
$ws.Range("E20").Value = '  -1.83%  '
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.707'
$ws.Range("E22").Value = '  -3.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.71'
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.60'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.32'
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.87'
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.89'
$ws.Range("E30").Value = '  -3.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("E31").Value = '  -2.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.60'
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("E35").Value = '  -4.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.90'
$ws.Range("E36").Value = '  -1.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.70'
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0393'
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.76'
$ws.Range("E40").Value = '  +3.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '409.87'
$ws.Range("E41").Value = '  -3.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.16'
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("E43").Value = '  -4.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.833.73'
$ws.Range("E44").Value = '  -8.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.258'
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.18'
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '36.17'
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.998'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.27'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.87'
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("E51").Value = '  -0.53%  '
